$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto symbol list values (price, volume%, hour, and row
# realignment due to inserted "One" row) as captured by the commit diff.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.39"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.68%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "21"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "28.23"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-4.35%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "21"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.66%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "21"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05700"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.56%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "21"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.633"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.83%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "21"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.204"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.40%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "21"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8513"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.71%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "21"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8958"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "3.17%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "21"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01006"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1,572.07%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "21"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1391"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.79%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "21"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07087"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.23%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "21"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03166"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "8.06%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "21"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09218"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.76%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "21"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001538"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.47%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "21"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005903"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.16%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "21"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.492"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.08%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "21"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-4.11%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "21"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3168"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.52%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "21"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03307"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.08%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "21"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1276"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.55%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "21"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.537"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.88%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "21"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04069"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.54%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "21"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.05%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "21"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001222"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.30%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "21"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-17.12%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "21"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001199"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.82%"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "21"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "21"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "21"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "21"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "21"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "21"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "21"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "21"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "21"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "21"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "21"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "21"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "21"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03792"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.06%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "21"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1064"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.62%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "21"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-35.21%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "21"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002199"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-9.38%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "21"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "22.43%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "21"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005269"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.49%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "21"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "21"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "62.31%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "21"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002269"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.46%"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "21"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "21"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "21"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "21"
